$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the trailing space in "Aantal " -> "Aantal"
$ws.Range("B2").Value = "Aantal"

# Set column widths as in the diff (AutoFit to best fit content)
$ws.Range("A1:C3").EntireColumn.AutoFit()

# Update selection to B2
$ws.Range("B2").Select()
